# Sierra Leone master data: update the lang_code column from "fra" (French,
# used for the Madagascar template) to "eng" (English, used for Sierra Leone).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A8").Value = "eng"

# The template picked up a trailing blank formatted row (row 9, column E)
# when it was adapted for the new country - replicate that extra row so the
# sheet's used range/dimension matches the updated template.
$ws.Range("E9").Font.Italic = $true
$ws.Range("E9").Font.Size = 10
$ws.Range("E9").Font.Color = 0
$ws.Range("E9").NumberFormat = "@"

# Header row height was also normalised to match the rest of the table.
$ws.Rows.Item(1).RowHeight = 15
